# Scheduled-runner refresh of the market-price / profit columns
# (currentAveragePrice, currentAveragePriceNQ, currentAveragePriceHQ,
# LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ -- columns H:N)
# across the per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).
# Values below are the freshly-fetched figures; a couple of rows also
# gain/lose a trailing HQ-profit cell (N) because that side of the
# recipe has no HQ variant for this pull.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2152.1538   # was 2169.359
$ws.Range("I19").Value = 3530.0557   # was 3531.6667
$ws.Range("J19").Value = 971.0952   # was 1001.6667
$ws.Range("K19").Value = 3530.0557   # was 3531.6667
$ws.Range("L19").Value = 971.0952   # was 1001.6667
$ws.Range("M19").Value = -3355.0557   # was -3356.6667
$ws.Range("N19").Value = -1321.0952   # was -1351.6667
$ws.Range("H80").Value = 10149.889   # was 16393.154
$ws.Range("I80").Value = 6404.0527   # was 6845.6665
$ws.Range("J80").Value = 19046.25   # was 37875
$ws.Range("K80").Value = 19212.1581   # was 20536.9995
$ws.Range("L80").Value = 57138.75   # was 113625
$ws.Range("M80").Value = -18214.1581   # was -19538.9995
$ws.Range("N80").Value = -59134.75   # was -115621
$ws.Range("H83").Value = 10149.889   # was 16393.154
$ws.Range("I83").Value = 6404.0527   # was 6845.6665
$ws.Range("J83").Value = 19046.25   # was 37875
$ws.Range("K83").Value = 57636.4743   # was 61610.9985
$ws.Range("L83").Value = 171416.25   # was 340875
$ws.Range("M83").Value = -52644.4743   # was -56618.9985
$ws.Range("N83").Value = -181400.25   # was -350859
$ws.Range("H103").Value = 1082.875   # was 1094.75
$ws.Range("I103").Value = 902   # was 1004
$ws.Range("J103").Value = 1143.1666   # was 1107.7142
$ws.Range("K103").Value = 2706   # was 3012
$ws.Range("L103").Value = 3429.4998   # was 3323.1426
$ws.Range("M103").Value = -2120   # was -2426
$ws.Range("N103").Value = -4601.4998   # was -4495.142599999999
$ws.Range("H113").Value = 2798.525   # was 2783.0715
$ws.Range("I113").Value = 1856.1666   # was 1914.4445
$ws.Range("J113").Value = 3569.5454   # was 3434.5417
$ws.Range("K113").Value = 1856.1666   # was 1914.4445
$ws.Range("L113").Value = 3569.5454   # was 3434.5417
$ws.Range("M113").Value = 1397.8334   # was 1339.5555
$ws.Range("N113").Value = -10077.5454   # was -9942.5417
$ws.Range("H132").Value = 5031.387   # was 4837.2144
$ws.Range("I132").Value = 1976.8096   # was 1165.6666
$ws.Range("K132").Value = 5930.4288   # was 3496.9998
$ws.Range("M132").Value = -3400.4288   # was -966.9998000000001
$ws.Range("H137").Value = 4984.091   # was 5597.0527
$ws.Range("I137").Value = 6794.7   # was 8268
$ws.Range("J137").Value = 3475.25   # was 3654.5454
$ws.Range("K137").Value = 20384.1   # was 24804
$ws.Range("L137").Value = 10425.75   # was 10963.6362
$ws.Range("M137").Value = -17834.1   # was -22254
$ws.Range("N137").Value = -15525.75   # was -16063.6362

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 833.0417   # was 1467.3636
$ws.Range("I97").Value = 833.0417   # was 1518
$ws.Range("J97").Value = 0   # was 961
$ws.Range("K97").Value = 833.0417   # was 1518
$ws.Range("L97").Value = 0   # was 961
$ws.Range("M97").Value = -337.0417   # was -1022
$ws.Range("N97").Value = $null   # was -1953

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1387   # was 1487.9231
$ws.Range("I20").Value = 1159.5834   # was 1239.375
$ws.Range("J20").Value = 1776.8572   # was 1885.6
$ws.Range("K20").Value = 1159.5834   # was 1239.375
$ws.Range("L20").Value = 1776.8572   # was 1885.6
$ws.Range("M20").Value = -912.5834   # was -992.375
$ws.Range("N20").Value = -2270.8572   # was -2379.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3801.658   # was 4450.2583
$ws.Range("I86").Value = 2809.4   # was 3721.7693
$ws.Range("J86").Value = 4904.1665   # was 4976.3887
$ws.Range("K86").Value = 2809.4   # was 3721.7693
$ws.Range("L86").Value = 4904.1665   # was 4976.3887
$ws.Range("M86").Value = -1686.4   # was -2598.7693
$ws.Range("N86").Value = -7150.1665   # was -7222.3887
$ws.Range("H89").Value = 3801.658   # was 4450.2583
$ws.Range("I89").Value = 2809.4   # was 3721.7693
$ws.Range("J89").Value = 4904.1665   # was 4976.3887
$ws.Range("K89").Value = 14047   # was 18608.8465
$ws.Range("L89").Value = 24520.8325   # was 24881.9435
$ws.Range("M89").Value = -8431   # was -12992.8465
$ws.Range("N89").Value = -35752.8325   # was -36113.9435
$ws.Range("H105").Value = 837.32355   # was 814.0769
$ws.Range("I105").Value = 793.5357   # was 790.3103599999999
$ws.Range("J105").Value = 1041.6666   # was 883
$ws.Range("K105").Value = 793.5357   # was 790.3103599999999
$ws.Range("L105").Value = 1041.6666   # was 883
$ws.Range("M105").Value = 953.4643   # was 956.6896400000001
$ws.Range("N105").Value = -4535.6666   # was -4377
$ws.Range("H122").Value = 1054.6154   # was 748.05
$ws.Range("I122").Value = 892.5   # was 762.41174
$ws.Range("J122").Value = 3000   # was 666.6667
$ws.Range("K122").Value = 2677.5   # was 2287.23522
$ws.Range("L122").Value = 9000   # was 2000.0001
$ws.Range("M122").Value = -227.5   # was 162.76478
$ws.Range("N122").Value = -13900   # was -6900.0001
$ws.Range("H132").Value = 2232.1562   # was 2829.95
$ws.Range("I132").Value = 1466.5883   # was 1874.4286
$ws.Range("J132").Value = 3099.8   # was 3344.4614
$ws.Range("K132").Value = 4399.7649   # was 5623.2858
$ws.Range("L132").Value = 9299.400000000001   # was 10033.3842
$ws.Range("M132").Value = -1869.7649   # was -3093.2858
$ws.Range("N132").Value = -14359.4   # was -15093.3842
$ws.Range("H134").Value = 1738   # was 2171.3333
$ws.Range("I134").Value = 1137.7142   # was 1218.091
$ws.Range("J134").Value = 2384.4614   # was 2977.923
$ws.Range("K134").Value = 3413.1426   # was 3654.273
$ws.Range("L134").Value = 7153.3842   # was 8933.769
$ws.Range("M134").Value = -878.1425999999997   # was -1119.273
$ws.Range("N134").Value = -12223.3842   # was -14003.769

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41.166668   # was 39.5
$ws.Range("I12").Value = 22.375   # was 22.222221
$ws.Range("J12").Value = 48   # was 46.260868
$ws.Range("K12").Value = 67.125   # was 66.666663
$ws.Range("L12").Value = 144   # was 138.782604
$ws.Range("M12").Value = 105.875   # was 106.333337
$ws.Range("N12").Value = -490   # was -484.782604
$ws.Range("H69").Value = 2271.2104   # was 2376.0588
$ws.Range("J69").Value = 3226.6667   # was 3596
$ws.Range("L69").Value = 9680.000100000001   # was 10788
$ws.Range("N69").Value = -11302.0001   # was -12410
$ws.Range("H72").Value = 2271.2104   # was 2376.0588
$ws.Range("J72").Value = 3226.6667   # was 3596
$ws.Range("L72").Value = 29040.0003   # was 32364
$ws.Range("N72").Value = -37152.0003   # was -40476
$ws.Range("H129").Value = 1354.3529   # was 1300.7949
$ws.Range("I129").Value = 645.5   # was 685.1539
$ws.Range("J129").Value = 1850.55   # was 1608.6154
$ws.Range("K129").Value = 1936.5   # was 2055.4617
$ws.Range("L129").Value = 5551.65   # was 4825.8462
$ws.Range("M129").Value = 3063.5   # was 2944.5383
$ws.Range("N129").Value = -15551.65   # was -14825.8462
$ws.Range("H131").Value = 2098.1167   # was 2133.8386
$ws.Range("J131").Value = 1624.695   # was 1676.5245
$ws.Range("L131").Value = 4874.085   # was 5029.5735
$ws.Range("N131").Value = -14954.085   # was -15109.5735
$ws.Range("H132").Value = 4439.6   # was 4320.1
$ws.Range("I132").Value = 4583.3335   # was 4128.5713
$ws.Range("J132").Value = 4224   # was 4767
$ws.Range("K132").Value = 41250.0015   # was 37157.14169999999
$ws.Range("L132").Value = 38016   # was 42903
$ws.Range("M132").Value = -38720.0015   # was -34627.14169999999
$ws.Range("N132").Value = -43076   # was -47963
$ws.Range("H134").Value = 5110.846   # was 3723.7222
$ws.Range("I134").Value = 4703.4165   # was 3968.4666
$ws.Range("J134").Value = 10000   # was 2500
$ws.Range("K134").Value = 14110.2495   # was 11905.3998
$ws.Range("L134").Value = 30000   # was 7500
$ws.Range("M134").Value = -9040.249500000002   # was -6835.399800000001
$ws.Range("N134").Value = -40140   # was -17640
$ws.Range("H137").Value = 5733.4287   # was 4671.2163
$ws.Range("I137").Value = 826.8   # was 785.4
$ws.Range("J137").Value = 18000   # was 12766.667
$ws.Range("K137").Value = 2480.4   # was 2356.2
$ws.Range("L137").Value = 54000   # was 38300.001
$ws.Range("M137").Value = 2619.6   # was 2743.8
$ws.Range("N137").Value = -64200   # was -48500.001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3896.3809   # was 4541.294
$ws.Range("I102").Value = 1835.4166   # was 2189.111
$ws.Range("J102").Value = 6644.3335   # was 7187.5
$ws.Range("K102").Value = 1835.4166   # was 2189.111
$ws.Range("L102").Value = 6644.3335   # was 7187.5
$ws.Range("M102").Value = -213.4166   # was -567.1109999999999
$ws.Range("N102").Value = -9888.333500000001   # was -10431.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1363.9032   # was 1502.5927
$ws.Range("I46").Value = 1168.4736   # was 1216.1111
$ws.Range("J46").Value = 1673.3334   # was 2075.5557
$ws.Range("K46").Value = 1168.4736   # was 1216.1111
$ws.Range("L46").Value = 1673.3334   # was 2075.5557
$ws.Range("M46").Value = -980.4736   # was -1028.1111
$ws.Range("N46").Value = -2049.3334   # was -2451.5557
$ws.Range("H122").Value = 1954.4546   # was 2501
$ws.Range("I122").Value = 2034   # was 2501
$ws.Range("J122").Value = 1859   # was 0
$ws.Range("K122").Value = 6102   # was 7503
$ws.Range("L122").Value = 5577   # was 0
$ws.Range("M122").Value = -3652   # was -5053
$ws.Range("N122").Value = -10477   # newly populated (cell previously empty)
